# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only the "K" column (column G) values change for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 2
